{"js": "const replacements = [\n  [\"2024-04-20 Saturday\", \"2024-04-21 Sunday\"],\n  [\"427\u00d76=2562\", \"492\u00d79=4428\"],\n  [\"602\u00d74=2408\", \"246\u00d78=1968\"],\n  [\"650\u00d75=3250\", \"829\u00d78=6632\"],\n  [\"402\u00d77=2814\", \"282\u00d76=1692\"],\n  [\"399\u00d73=1197\", \"342\u00d75=1710\"],\n  [\"618\u00d72=1236\", \"909\u00d74=3636\"],\n  [\"896\u00d73=2688\", \"648\u00d78=5184\"],\n  [\"658\u00d76=3948\", \"285\u00d72=570\"],\n  [\"654\u00d79=5886\", \"623\u00d79=5607\"],\n  [\"572\u00d74=2288\", \"349\u00d74=1396\"],\n  [\"246\u00d75=1230\", \"233\u00d72=466\"],\n  [\"980\u00d73=2940\", \"622\u00d72=1244\"],\n  [\"249\u00d76=1494\", \"660\u00d76=3960\"],\n  [\"752\u00d74=3008\", \"374\u00d75=1870\"],\n  [\"439\u00d76=2634\", \"364\u00d76=2184\"],\n  [\"970\u00d79=8730\", \"529\u00d78=4232\"],\n  [\"323\u00d79=2907\", \"420\u00d73=1260\"],\n  [\"719\u00d73=2157\", \"932\u00d79=8388\"],\n  [\"278\u00d77=1946\", \"694\u00d79=6246\"],\n  [\"195\u00d78=1560\", \"736\u00d77=5152\"],\n  [\"662\u00d75=3310\", \"465\u00d73=1395\"],\n  [\"251\u00d73=753\", \"274\u00d72=548\"],\n  [\"577\u00d75=2885\", \"122\u00d78=976\"],\n  [\"997\u00d76=5982\", \"197\u00d76=1182\"],\n  [\"993\u00d73=2979\", \"529\u00d79=4761\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $found = $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n  if (-not $found) {\n    throw \"Could not find text: $findText\"\n  }\n}\n\nReplace-Text \"2024-04-20 Saturday\" \"2024-04-21 Sunday\"\nReplace-Text \"427\u00d76=2562\" \"492\u00d79=4428\"\nReplace-Text \"602\u00d74=2408\" \"246\u00d78=1968\"\nReplace-Text \"650\u00d75=3250\" \"829\u00d78=6632\"\nReplace-Text \"402\u00d77=2814\" \"282\u00d76=1692\"\nReplace-Text \"399\u00d73=1197\" \"342\u00d75=1710\"\nReplace-Text \"618\u00d72=1236\" \"909\u00d74=3636\"\nReplace-Text \"896\u00d73=2688\" \"648\u00d78=5184\"\nReplace-Text \"658\u00d76=3948\" \"285\u00d72=570\"\nReplace-Text \"654\u00d79=5886\" \"623\u00d79=5607\"\nReplace-Text \"572\u00d74=2288\" \"349\u00d74=1396\"\nReplace-Text \"246\u00d75=1230\" \"233\u00d72=466\"\nReplace-Text \"980\u00d73=2940\" \"622\u00d72=1244\"\nReplace-Text \"249\u00d76=1494\" \"660\u00d76=3960\"\nReplace-Text \"752\u00d74=3008\" \"374\u00d75=1870\"\nReplace-Text \"439\u00d76=2634\" \"364\u00d76=2184\"\nReplace-Text \"970\u00d79=8730\" \"529\u00d78=4232\"\nReplace-Text \"323\u00d79=2907\" \"420\u00d73=1260\"\nReplace-Text \"719\u00d73=2157\" \"932\u00d79=8388\"\nReplace-Text \"278\u00d77=1946\" \"694\u00d79=6246\"\nReplace-Text \"195\u00d78=1560\" \"736\u00d77=5152\"\nReplace-Text \"662\u00d75=3310\" \"465\u00d73=1395\"\nReplace-Text \"251\u00d73=753\" \"274\u00d72=548\"\nReplace-Text \"577\u00d75=2885\" \"122\u00d78=976\"\nReplace-Text \"997\u00d76=5982\" \"197\u00d76=1182\"\nReplace-Text \"993\u00d73=2979\" \"529\u00d79=4761\"\n"}
